$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table (rows 54-58) by copying the formatting of the last existing
# data row (53) down, then writing in the new values - this keeps per-column
# styles (s="1" on A, s="2" on B/C date cells, s="3" on D/G/H/I/J/K) identical
# to the rest of the table, matching how the sheet was originally built.
$ws.Range("A53:L53").Copy() | Out-Null
$ws.Range("A54:L58").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 54
$ws.Range("A54").Value2 = 52
$ws.Range("B54").Value2 = 44896.3125
$ws.Range("C54").Value2 = 44896.33333333334
$ws.Range("D54").Value2 = 0.02083333333333333
$ws.Range("E54").Value2 = "t_ar"
$ws.Range("F54").Value2 = 6.764646464646465
$ws.Range("G54").Value2 = 0
$ws.Range("H54").Value2 = 0.02083333333333333
$ws.Range("I54").Value2 = 0.9791666666666666
$ws.Range("J54").Value2 = 5.75
$ws.Range("K54").Value2 = 77.04166666666667
$ws.Range("L54").Value2 = 1849

# Row 55
$ws.Range("A55").Value2 = 53
$ws.Range("B55").Value2 = 44897.3125
$ws.Range("C55").Value2 = 44897.33333333334
$ws.Range("D55").Value2 = 0.02083333333333333
$ws.Range("E55").Value2 = "t_ar"
$ws.Range("F55").Value2 = 6.764646464646465
$ws.Range("G55").Value2 = 0
$ws.Range("H55").Value2 = 0.02083333333333333
$ws.Range("I55").Value2 = 3.020833333333333
$ws.Range("J55").Value2 = 5.770833333333333
$ws.Range("K55").Value2 = 80.0625
$ws.Range("L55").Value2 = 1921.5

# Row 56
$ws.Range("A56").Value2 = 54
$ws.Range("B56").Value2 = 44900.35416666666
$ws.Range("C56").Value2 = 44900.375
$ws.Range("D56").Value2 = 0.02083333333333333
$ws.Range("E56").Value2 = "t_ar"
$ws.Range("F56").Value2 = 6.764646464646465
$ws.Range("G56").Value2 = 0
$ws.Range("H56").Value2 = 0.02083333333333333
$ws.Range("I56").Value2 = 0.4375
$ws.Range("J56").Value2 = 5.791666666666667
$ws.Range("K56").Value2 = 80.5
$ws.Range("L56").Value2 = 1932

# Row 57
$ws.Range("A57").Value2 = 55
$ws.Range("B57").Value2 = 44900.8125
$ws.Range("C57").Value2 = 44900.83333333334
$ws.Range("D57").Value2 = 0.02083333333333333
$ws.Range("E57").Value2 = "t_ar"
$ws.Range("F57").Value2 = 6.764646464646465
$ws.Range("G57").Value2 = 0
$ws.Range("H57").Value2 = 0.02083333333333333
$ws.Range("I57").Value2 = 0.4791666666666667
$ws.Range("J57").Value2 = 5.8125
$ws.Range("K57").Value2 = 80.97916666666667
$ws.Range("L57").Value2 = 1943.5

# Row 58
$ws.Range("A58").Value2 = 56
$ws.Range("B58").Value2 = 44901.3125
$ws.Range("C58").Value2 = 44901.33333333334
$ws.Range("D58").Value2 = 0.02083333333333333
$ws.Range("E58").Value2 = "t_ar"
$ws.Range("F58").Value2 = 6.764646464646465
$ws.Range("G58").Value2 = 0
$ws.Range("H58").Value2 = 0.02083333333333333
$ws.Range("I58").Value2 = 7.5
$ws.Range("J58").Value2 = 5.833333333333333
$ws.Range("K58").Value2 = 88.47916666666667
$ws.Range("L58").Value2 = 2123.5

